$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of statistics data to append below the existing data (row 2)
$data = @(
    @(2, "TestName TestSurname", "TestParentName TestParentSurname", "085 555 5555", "Eerste"),
    @(3, "TestName2 TestSurname2", "TestParentName3 TestParentSurname3", "085 555 6666", "Tweede"),
    @(4, "TestName3 TestSurname3", "TestParentName3 TestParentSurname3", "085 555 6666", "Tweede"),
    @(5, "TestName TestSurname", "TestParentName TestParentSurname", "085 555 6666", "Tweede"),
    @(6, "TestName TestSurname", "TestParentName TestParent Surname", "085 555 6666", "Tweede")
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
